$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1369994
$ws.Range("C4").Value = 2356
$ws.Range("E4").Value = 1032801
$ws.Range("G4").Value = 61
$ws.Range("H4").Value = 80848

# Row 16 - India
$ws.Range("B16").Value = 68789
$ws.Range("C16").Value = 1628
$ws.Range("D16").Value = 21266
$ws.Range("E16").Value = 45300
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 2223

# Row 19 - Paises Bajos
$ws.Range("F19").Value = 498

# Row 44 - Dinamarca
$ws.Range("D44").Value = 8328
$ws.Range("E44").Value = 1652
$ws.Range("F44").Value = 43
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 533

# Row 79 - Bulgaria
$ws.Range("B79").Value = 1990
$ws.Range("C79").Value = 25
$ws.Range("E79").Value = 1436
$ws.Range("F79").Value = 50
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 93

# Row 116 - Kenia
$ws.Range("D116").Value = 251
$ws.Range("E116").Value = 416
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 33
